# Automatische test-sync: 2025-08-04 21:02:50
# Append a new test-mail row (row 29) to the "Logs" sheet and refresh the
# related conditional-formatting ranges + the "Dashboard" summary count.

$wb = $excel.ActiveWorkbook

# ---- 1. Logs sheet: append row 29 ----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A29").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("B29").Value = "mailmind.test@zohomail.eu"
$logs.Range("C29").Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D29").Value = "Planning / Afspraak"
$logs.Range("E29").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F29").Value = "2025-08-04 21:02:31"
$logs.Range("G29").Value = "Ja"
$logs.Range("H29").Value = "Ja"
$logs.Range("I29").Value = "Nee"
$logs.Range("J29").Value = "Nee"

# ---- 2. Extend the conditional formatting ranges from row 28 to row 29 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "28")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "29")
    $count = $oldRange.FormatConditions.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $oldRange.FormatConditions.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# ---- 3. Dashboard sheet: bump the "Planning / Afspraak" count ------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
